$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits, on its own, in the
#    empty paragraph right after the "Cristian, Artur, Si et. al" line.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# Helper: clears a just-created (InsertParagraphAfter) empty paragraph
# down to *no* run at all (matches the document's own convention for
# blank paragraphs, which carry no <w:r>). A placeholder character has
# to be written first -- clearing an already-empty paragraph mark is a
# no-op in this engine and leaves a phantom empty <w:r> behind.
# ------------------------------------------------------------------
function Clear-EmptyParagraph($para) {
    $para.Range.Text = "z"
    $rng = $para.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = ""
}

# ------------------------------------------------------------------
# 2) Append six new paragraphs after "... outside the scope of this
#    article." (currently the last paragraph of the reviewer-response
#    body, right before the trailing indented blank paragraph).
# ------------------------------------------------------------------

$lastReply = $d.Paragraphs.Item(29)
$lastReply.Range.InsertParagraphAfter()

# --- paragraph: empty -------------------------------------------------
$pA = $d.Paragraphs.Item(30)
Clear-EmptyParagraph $pA

# --- paragraph: "We add the following sentences ..." ------------------
$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Item(31)
$pB.Range.Text = "We add the following sentences in order to address the comment (right at the end of line 94):"

# --- paragraph: empty ---------------------------------------------------
$pB.Range.InsertParagraphAfter()
$pC = $d.Paragraphs.Item(32)
Clear-EmptyParagraph $pC

# --- paragraph: quoted addition (bold, green, first-line indent) -------
$pC.Range.InsertParagraphAfter()
$pD = $d.Paragraphs.Item(33)
$pD.Format.FirstLineIndent = 36
$pD.Range.Text = "“… We choose a relatively simple impulse response"
$pD.Range.Font.Bold = $true
$pD.Range.Font.Color = 5287936

# Text laid down above is: "<U+201C><U+2026> We choose a relatively..."
#                            0        1      2 3
# Run boundaries wanted:  [0,1) "\u201c" | [1,3) "\u2026 " | [3,end) "We..."
$pDrange = $d.Paragraphs.Item(33).Range
$splitAt1 = $pDrange.Start + 1
$bmrng1 = $d.Range($splitAt1, $splitAt1)
$d.Bookmarks.Add("TmpSplitD1", $bmrng1) | Out-Null
$d.Bookmarks.Item("TmpSplitD1").Delete()

$pDrange2 = $d.Paragraphs.Item(33).Range
$splitAt2 = $pDrange2.Start + 3
$gobackrng = $d.Range($splitAt2, $splitAt2)
$d.Bookmarks.Add("_GoBack", $gobackrng) | Out-Null

# --- paragraph: "function in order to study ..." ------------------------
$pD2 = $d.Paragraphs.Item(33)
$pD2.Range.InsertParagraphAfter()
$pE = $d.Paragraphs.Item(34)
$pE.Format.FirstLineIndent = 0
$pE.Range.Text = " function in order to study the general behavior of the time resolution, "
$pE.Range.Font.Bold = $true
$pE.Range.Font.Color = 5287936

# --- paragraph: "while detailed studies ..." + closing quote ------------
$pE.Range.InsertParagraphAfter()
$pF = $d.Paragraphs.Item(35)
$pF.Format.FirstLineIndent = 0
$pF.Range.Text = " while detailed studies of more realistic circuit implementations and correlated effects are left for future studies.” "
$pF.Range.Font.Bold = $true
$pF.Range.Font.Color = 5287936

$pFrange = $d.Paragraphs.Item(35).Range
$splitAt3 = $pFrange.End - 3
$bmrng3 = $d.Range($splitAt3, $splitAt3)
$d.Bookmarks.Add("TmpSplitF1", $bmrng3) | Out-Null
$d.Bookmarks.Item("TmpSplitF1").Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
